$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the source matrix A (B2:D4)
$ws.Range("B2").Value = 7
$ws.Range("C2").Value = -6
$ws.Range("D2").Value = 6

$ws.Range("B3").Value = 2
$ws.Range("C3").Value = -1
$ws.Range("D3").Value = 2

$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = -1

# Update manually entered eigenvalues lambda2 (C20) and lambda3 (C21)
$ws.Range("C20").Value = -3
$ws.Range("C21").Value = 7

# Fix lambda3 reference cell (C41) to be dynamic: reference C21 instead of a static number
$ws.Range("C41").Formula = "=C21"

# Update the view so the visible window matches the target state
$ws.Range("C42").Select()
